# Add ability to schedule a message to be sent in the future
# -> add a 4th "SendTime(optional)" column to the bulk-sms template,
#    with a sample datetime value on the 2nd example row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + sample value for the scheduled send-time column.
$ws.Range("D1").Value = "SendTime(optional)"
$ws.Range("D3").Value = "2023-11-11T02:10:01"

# Re-balance the column widths now that there are 4 columns instead of 3.
$ws.Columns.Item(3).ColumnWidth = 33.6328125
$ws.Columns.Item(4).ColumnWidth = 20.36328125

# Leave the selection where the author left it while editing the sheet.
$ws.Range("G8").Select() | Out-Null
